$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add "U" sex value in column D for rows 2-6, and clear the now-unused
# sire/dam numeric columns (F and G).
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 4).Value = "U"
    $ws.Cells.Item($r, 6).ClearContents()
    $ws.Cells.Item($r, 7).ClearContents()
}

# Update the selected cell to match the saved view state.
$ws.Range("F16").Select()
